$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 updates
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5

# Row 4 updates
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 3
$ws.Range("AC4").Value = 12
$ws.Range("AD4").Value = 7
$ws.Range("AI4").Value = 17
$ws.Range("AL4").Value = 26
$ws.Range("AQ4").Value = 34
$ws.Range("AS4").Value = 126
$ws.Range("AT4").Value = 3
$ws.Range("AU4").Value = 7.5
